# Daily attendance processing - 2026-01-25 07:13:53
#
# The "Recorded By" column (G) lists the users/services that recorded a
# session as a comma-separated string. For rows where the automated
# "system"/"System" account was recorded FIRST in that list, move it to the
# end so the human editor(s) are listed first, e.g.:
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "system, backup@backdoor.com, System"     -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text -eq "System, dnasr281@gmail.com" -or $text -eq "system, backup@backdoor.com, System") {
        $parts = $text -split ", "
        $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
        $cell.Value = $rotated -join ", "
    }
}
